# 5.5.2.xlsx: add the 2020 data point (Q4/Q5) and update the view state
# (scrolled/selected cell) to match the published workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# --- New "2020" column -------------------------------------------------
# Q4 header (year) and Q5 value, copying the formatting that P4/P5 (the
# previous "2019" column) already carry so the new cells pick up the same
# style indices instead of a bare default style.
$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)
$ws.Range("Q4").Value = 2020

$ws.Range("P5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)
$ws.Range("Q5").Value = 47.4

# --- View state ----------------------------------------------------------
# Scroll the viewport so column C is left-most and leave the cursor parked
# on Q9 (one row below the new data), same as the authored workbook.
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("Q9").Select()

"done"
